# Add new Broadcom BCM57414 NIC card compatibility row to the
# "openEuler22.03-LTS两类平台板卡兼容性" worksheet (row 46), and extend the
# autofilter / filter-database range to cover it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Seed row 46 from row 45 so it inherits the same base formatting
#     (borders, fonts, number formats) as the rest of the data table. ---
$ws.Range("A45:Q45").Copy($ws.Range("A46:Q46"))

# --- A few source cells already carry the exact cell styles the new
#     row needs for specific columns; grab their formats. ---
$ws.Range("C42").Copy()
$ws.Range("C46").PasteSpecial(-4122)

$ws.Range("O21").Copy()
$ws.Range("K46:O46").PasteSpecial(-4122)

$ws.Range("L42").Copy()
$ws.Range("L46").PasteSpecial(-4122)

# --- Values ---
$ws.Cells.Item(46, 1).Value = "14e4"                          # A vendorID
$ws.Cells.Item(46, 2).Value = "16d7"                          # B deviceID
$ws.Cells.Item(46, 3).Value = "14e4"                          # C svID
$ws.Cells.Item(46, 4).Value = "1402"                          # D ssID
$ws.Cells.Item(46, 5).Value = "aarch64"                       # E architecture
$ws.Cells.Item(46, 6).Value = "openEuler 22.03 LTS"           # F os
$ws.Cells.Item(46, 7).Value = "bnxt_en"                       # G driverName
$ws.Cells.Item(46, 9).Value = "NIC"                           # I type

# J (Date) looks like a date to the parser ("2022.05.26"), so park it in
# a Text-formatted cell first, then restore the General-format style
# that the rest of the Date column uses (style is re-pasted, value stays).
$ws.Cells.Item(46, 10).NumberFormat = "@"
$ws.Cells.Item(46, 10).Value = "2022.05.26"                   # J Date
$ws.Range("J45").Copy()
$ws.Range("J46").PasteSpecial(-4122)

$ws.Cells.Item(46, 11).Value = "5207352D392184F5F7FBE52CE0A5C064DA1687B2"  # K sha256
$ws.Cells.Item(46, 12).Value = "406K"                         # L driverSize
$ws.Cells.Item(46, 13).Value = "Broadcom"                     # M chipVendor
$ws.Cells.Item(46, 14).Value = "BCM957414A4142CC_08"          # N boardModel
$ws.Cells.Item(46, 15).Value = "BCM57414"                     # O chipModel
$ws.Cells.Item(46, 16).Value = "06310148"                     # P item
$ws.Cells.Item(46, 17).Value = "inbox"                        # Q downloadLink

# --- Rich-text runs: first letter default, remainder explicitly
#     re-applied (mirrors how the source workbook stores them). ---
$kChars = $ws.Range("K46").Characters(2, 39)
$kChars.Font.Bold = $false

$mChars = $ws.Range("M46").Characters(2, 7)
$mChars.Font.Bold = $false

$oChars = $ws.Range("O46").Characters(2, 7)
$oChars.Font.Bold = $false

# --- Q46 in the source sheet drops the top border (it's now the final
#     row), matching the generated style used for the table edge. ---
$ws.Range("Q46").Borders.Item(8).LineStyle = -4142

# --- Extend the AutoFilter + hidden _FilterDatabase name to the new
#     used range, A1:U46. ---
$ws.AutoFilterMode = $false
$ws.Range("A1:U46").AutoFilter()
$wb.Names.Item("_xlnm._FilterDatabase").RefersTo = "='openEuler22.03-LTS两类平台板卡兼容性'!`$A`$1:`$U`$46"
